$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.28768789768219
$ws.Range("B1").Value = 1.421594977378845
$ws.Range("C1").Value = 1.695623397827148
$ws.Range("D1").Value = 3.041944980621338
$ws.Range("E1").Value = 4.258532524108887
